# Update the "Riders" (C) and "Average" (D) columns on the Ridership sheet
# with new Madigan bike hours data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Ridership")

# Row -> (Riders, Average)
$ws.Range("C2").Value = 289
$ws.Range("D2").Value = 267.57

$ws.Range("C3").Value = 217
$ws.Range("D3").Value = 228.86

$ws.Range("C4").Value = 190
$ws.Range("D4").Value = 197.14

$ws.Range("C5").Value = 241
$ws.Range("D5").Value = 228.83

$ws.Range("C6").Value = 332
$ws.Range("D6").Value = 228.29

$ws.Range("C7").Value = 181
$ws.Range("D7").Value = 132.5

$ws.Range("C8").Value = 145
$ws.Range("D8").Value = 105.5

$wb.Save()
